$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.103.53'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '3.245.09'
$ws.Range("E3").Value = '  +3.01%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''592.57'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '''140.45'
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.245.54'
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").Value = '''0.519'
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").Value = '''0.147'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = '''5.31'
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").Value = '''0.463'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '''0.0000247'
$ws.Range("E13").Value = '  -2.91%  '
$ws.Range("D14").Value = '''34.40'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '3.778.58'
$ws.Range("E15").Value = '  +3.01%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '3.244.91'
$ws.Range("E17").Value = '  +3.09%  '
$ws.Range("D18").Value = '63.180.76'
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").Value = '''6.74'
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("D20").Value = '''475.07'
$ws.Range("E20").Value = '  -2.84%  '
$ws.Range("D21").Value = '''14.10'
$ws.Range("E21").Value = '  -3.94%  '
$ws.Range("D22").Value = '''0.719'
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").Value = '''7.90'
$ws.Range("E23").Value = '  +3.29%  '
$ws.Range("D24").Value = '''83.70'
$ws.Range("E24").Value = '  -4.58%  '
$ws.Range("D25").Value = '''13.24'
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("E28").Value = '  +4.01%  '
$ws.Range("D29").Value = '''8.06'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").Value = '''2.12'
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("D31").Value = '''27.46'
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("E33").Value = '  -2.95%  '
$ws.Range("D34").Value = '''2.53'
$ws.Range("E34").Value = '  -4.18%  '
$ws.Range("D35").Value = '''1.10'
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("D36").Value = '''5.86'
$ws.Range("E36").Value = '  -2.82%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").Value = '0.0₃0714'
$ws.Range("E38").Value = '  -4.39%  '
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("D40").Value = '''419.68'
$ws.Range("E40").Value = '  -3.94%  '
$ws.Range("D41").Value = '2.989.19'
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("D42").Value = '''8.37'
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("E43").Value = '  -7.80%  '
$ws.Range("D44").Value = '''0.111'
$ws.Range("E44").Value = '  -7.35%  '
$ws.Range("D45").Value = '''0.266'
$ws.Range("E45").Value = '  +2.83%  '
$ws.Range("D46").Value = '''2.15'
$ws.Range("E46").Value = '  -2.22%  '
$ws.Range("D48").Value = '''25.79'
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.114'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''2.29'
$ws.Range("E50").Value = '  -4.66%  '
$ws.Range("D51").Value = '''119.15'
$ws.Range("E51").Value = '  -1.07%  '
